$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.617.75'
$ws.Range('E2').Value = '  -4.42%  '
$ws.Range('D3').Value = '3.290.51'
$ws.Range('E3').Value = '  -6.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.44'
$ws.Range('E5').Value = '  -4.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.18'
$ws.Range('E6').Value = '  -11.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.281.50'
$ws.Range('E8').Value = '  -7.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  -10.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.174'
$ws.Range('E10').Value = '  -12.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.77'
$ws.Range('E11').Value = '  -6.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.515'
$ws.Range('E12').Value = '  -12.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.94'
$ws.Range('E13').Value = '  -15.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000247'
$ws.Range('E14').Value = '  -10.44%  '
$ws.Range('D15').Value = '3.819.58'
$ws.Range('D16').Value = '67.614.20'
$ws.Range('E16').Value = '  -4.61%  '
$ws.Range('D17').Value = '3.287.65'
$ws.Range('E17').Value = '  -6.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.32'
$ws.Range('E18').Value = '  -13.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '538.38'
$ws.Range('E19').Value = '  -11.35%  '
$ws.Range('E20').Value = '  -6.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.20'
$ws.Range('E21').Value = '  -14.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.767'
$ws.Range('E22').Value = '  -12.96%  '
$ws.Range('E23').Value = '  -12.80%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.74'
$ws.Range('E24').Value = '  -11.97%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.18'
$ws.Range('E25').Value = '  -11.63%  '
$ws.Range('E26').Value = '  -11.58%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.20'
$ws.Range('E28').Value = '  -9.77%  '
$ws.Range('E29').Value = '  -15.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.56'
$ws.Range('E30').Value = '  -12.26%  '
$ws.Range('E31').Value = '  -9.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.16'
$ws.Range('E32').Value = '  -10.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '549.75'
$ws.Range('E33').Value = '  -10.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.68'
$ws.Range('E34').Value = '  -17.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.83'
$ws.Range('E35').Value = '  -14.32%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0466'
$ws.Range('E37').Value = '  -5.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.68'
$ws.Range('E38').Value = '  -5.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0870'
$ws.Range('E39').Value = '  -12.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.19'
$ws.Range('E40').Value = '  -15.33%  '
$ws.Range('E41').Value = '  -9.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.79'
$ws.Range('E42').Value = '  -17.79%  '
$ws.Range('D43').Value = '2.954.56'
$ws.Range('E43').Value = '  -11.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.274'
$ws.Range('E44').Value = '  -11.81%  '
$ws.Range('D45').Value = '0.0₃0604'
$ws.Range('E45').Value = '  -16.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.23'
$ws.Range('E46').Value = '  -10.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.09'
$ws.Range('E47').Value = '  -15.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.38'
$ws.Range('E48').Value = '  -17.82%  '
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '126.75'
$ws.Range('E50').Value = '  -5.34%  '
$ws.Range('E51').Value = '  -11.95%  '
